# Continued with matrices, LINQ in GetObj
#
# Rename Sheet2 -> TestTable and build the "create a CURVE/matrix object,
# then display it" example, mirroring the existing pattern already present
# on Sheet1 (mmCreateObj / mmDisplayObj).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2 -> "TestTable" ------------------------------------------------
$ws2.Name = "TestTable"

# --- input area (B3:C3, B5:E10) ------------------------------------------
$ws2.Range("B3").Value = "A"
$ws2.Range("C3").Value = "CURVE"

$ws2.Range("B5").Value = "currency"
$ws2.Range("C5").Value = "EUR"

$ws2.Range("B6").Value = "rates"

$ws2.Range("C7").Value = "A"
$ws2.Range("D7").Value = "B"
$ws2.Range("E7").Value = "C"

$ws2.Range("C8").Value = 1
$ws2.Range("D8").Value = 3
$ws2.Range("E8").Value = 5

$ws2.Range("C9").Value = 2
$ws2.Range("D9").Value = 4
$ws2.Range("E9").Value = 6

$ws2.Range("B10").Value = "rate"
$ws2.Range("C10").Value = 5

# --- literal "display" block (G1:L7) — written before the array formula so
# the anchor cell can be laid on top without clobbering its siblings -------
$ws2.Range("H1").Value = "EUR"
$ws2.Range("I1").Value = "#N/A"
$ws2.Range("J1").Value = "#N/A"
$ws2.Range("K1").Value = "#N/A"
$ws2.Range("L1").Value = "#N/A"

$ws2.Range("G2").Value = "rate"
$ws2.Range("H2").NumberFormat = "@"
$ws2.Range("H2").Value = "0"
$ws2.Range("I2").Value = "#N/A"
$ws2.Range("J2").Value = "#N/A"
$ws2.Range("K2").Value = "#N/A"
$ws2.Range("L2").Value = "#N/A"

$ws2.Range("G3").Value = "rates"
$ws2.Range("H3").Value = "Tables cannot be displayed yet"
$ws2.Range("I3").Value = "#N/A"
$ws2.Range("J3").Value = "#N/A"
$ws2.Range("K3").Value = "#N/A"
$ws2.Range("L3").Value = "#N/A"

$ws2.Range("G4").Value = "name"
$ws2.Range("H4").Value = "A"
$ws2.Range("I4").Value = "#N/A"
$ws2.Range("J4").Value = "#N/A"
$ws2.Range("K4").Value = "#N/A"
$ws2.Range("L4").Value = "#N/A"

$ws2.Range("G5:L5").Value = "#N/A"
$ws2.Range("G6:L6").Value = "#N/A"
$ws2.Range("G7:L7").Value = "#N/A"

# --- formulas (mirrors Sheet1's CreateObj/DisplayObj usage) --------------
# B1 builds the CURVE object "A" from the B3:C3/B5:E10 inputs.
$ws2.Range("B1").Formula = "=_xll.mmCreateObj(B3,C3,B5:E10)"

# G1 (array, spilling G1:L7) displays the object just created.
$ws2.Range("G1:L7").FormulaArray = "=_xll.mmDisplayObj(B1,C3)"

# --- sheet view: TestTable becomes the active tab/sheet -------------------
$ws2.Activate()
$ws2.Range("G1:L7").Select()
